# Refresh the cryptocurrency price/volume snapshot (D: Price, E: Volume(1h)).
# D-column values that look like plain decimals are written with a leading
# apostrophe so Excel keeps them as text (matching the source data's text type)
# instead of silently re-parsing them into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.906.34'
$ws.Range("E2").Value = '  -0.21%  '

$ws.Range("D3").Value = '2.977.84'
$ws.Range("E3").Value = '  -1.08%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '''498.79'
$ws.Range("E5").Value = '  -2.82%  '

$ws.Range("D6").Value = '''137.10'
$ws.Range("E6").Value = '  -1.81%  '

$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").Value = '''0.428'
$ws.Range("E8").Value = '  -2.04%  '

$ws.Range("E9").Value = '  -1.79%  '

$ws.Range("E10").Value = '  -0.74%  '

$ws.Range("D11").Value = '''0.357'
$ws.Range("E11").Value = '  -0.47%  '

$ws.Range("D12").Value = '3.494.25'
$ws.Range("E12").Value = '  -0.88%  '

$ws.Range("E13").Value = '  -1.37%  '

$ws.Range("D14").Value = '''25.79'
$ws.Range("E14").Value = '  -0.92%  '

$ws.Range("E15").Value = '  +0.35%  '

$ws.Range("D16").Value = '56.940.19'
$ws.Range("E16").Value = '  -0.15%  '

$ws.Range("E17").Value = '  +1.79%  '

$ws.Range("D18").Value = '2.972.82'
$ws.Range("E18").Value = '  -1.21%  '

$ws.Range("D19").Value = '''12.61'
$ws.Range("E19").Value = '  +0.19%  '

$ws.Range("D20").Value = '''7.79'
$ws.Range("E20").Value = '  -1.13%  '

$ws.Range("D21").Value = '''320.17'
$ws.Range("E21").Value = '  -2.38%  '

$ws.Range("E22").Value = '  -0.09%  '

$ws.Range("E23").Value = '  -0.93%  '

$ws.Range("D24").Value = '''0.487'
$ws.Range("E24").Value = '  -0.33%  '

$ws.Range("D25").Value = '''63.56'
$ws.Range("E25").Value = '  -0.26%  '

$ws.Range("E26").Value = '  +0.20%  '

$ws.Range("D27").Value = '''0.162'
$ws.Range("E27").Value = '  -5.95%  '

$ws.Range("D28").Value = '0.0₃0891'
$ws.Range("E28").Value = '  -2.94%  '

$ws.Range("D29").Value = '''6.58'
$ws.Range("E29").Value = '  -2.16%  '

$ws.Range("D30").Value = '''7.13'
$ws.Range("E30").Value = '  +0.68%  '

$ws.Range("E31").Value = '  -2.61%  '

$ws.Range("D32").Value = '''1.16'
$ws.Range("E32").Value = '  -6.67%  '

$ws.Range("D33").Value = '''20.15'
$ws.Range("E33").Value = '  -2.08%  '

$ws.Range("D34").Value = '''152.25'
$ws.Range("E34").Value = '  -2.51%  '

$ws.Range("E35").Value = '  +0.46%  '

$ws.Range("D36").Value = '''5.75'
$ws.Range("E36").Value = '  +0.33%  '

$ws.Range("D37").Value = '''1.25'
$ws.Range("E37").Value = '  -2.75%  '

$ws.Range("D38").Value = '''24.01'
$ws.Range("E38").Value = '  +0.27%  '

$ws.Range("D39").Value = '''0.0664'
$ws.Range("E39").Value = '  -2.53%  '

$ws.Range("D40").Value = '3.008.44'
$ws.Range("E40").Value = '  -1.13%  '

$ws.Range("D41").Value = '''37.46'

$ws.Range("D42").Value = '''0.999'
$ws.Range("E42").Value = '  -0.11%  '

$ws.Range("E43").Value = '  +1.21%  '

$ws.Range("D44").Value = '''0.641'
$ws.Range("E44").Value = '  -1.20%  '

$ws.Range("D45").Value = '2.199.03'
$ws.Range("E45").Value = '  -4.42%  '

$ws.Range("E46").Value = '  -3.04%  '

$ws.Range("D47").Value = '''0.949'
$ws.Range("E47").Value = '  -5.90%  '

$ws.Range("D48").Value = '''5.94'
$ws.Range("E48").Value = '  +0.64%  '

$ws.Range("D49").Value = '''0.0234'
$ws.Range("E49").Value = '  -2.29%  '

$ws.Range("D50").Value = '''19.13'
$ws.Range("E50").Value = '  -1.92%  '

$ws.Range("E51").Value = '  -8.34%  '
